# GDE-9324 - added other deals
# Adds 4 new UAT deal rows (rowid 2-5) to the "Clients" sheet, mirroring
# the formatting of the existing data row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clients")

# Copy the formatting (styles) of the existing data row down into the four
# new rows before filling in values, so the new cells pick up the same
# quote-prefixed "rowid" style and the hyperlink-like "Path" style used on
# row 2.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D6").PasteSpecial(-4122)

# Fill in the "Filename" column (D) for each new deal.
$ws.Cells.Item(3, 4).Value = "Deal_CH_EDU_BILAT.xlsx"
$ws.Cells.Item(4, 4).Value = "Deal_LBT_BILAT.xlsx"
$ws.Cells.Item(5, 4).Value = "Deal_New_Life_BILAT.xlsx"
$ws.Cells.Item(6, 4).Value = "Deal_PIM_Future_BILAT.xlsx"

# Fill in the "UAT_Client" column (B) for each new deal.
$ws.Cells.Item(4, 2).Value = "LBT_BILAT"
$ws.Cells.Item(5, 2).Value = "New_Life_BILAT"
$ws.Cells.Item(6, 2).Value = "PIM_Future_BILAT"
$ws.Cells.Item(3, 2).Value = "CH_EDU_BILAT"

# Fill in the "rowid" column (A) as text (quote-prefixed), continuing on
# from the existing rowid of 1.
$ws.Cells.Item(3, 1).Value = "'2"
$ws.Cells.Item(4, 1).Value = "'3"
$ws.Cells.Item(5, 1).Value = "'4"
$ws.Cells.Item(6, 1).Value = "'5"

# Fill in the "Path" column (C) - same path used by every deal.
$ws.Cells.Item(3, 3).Value = "\DataSet\NewUATDeals_DataSet\"
$ws.Cells.Item(4, 3).Value = "\DataSet\NewUATDeals_DataSet\"
$ws.Cells.Item(5, 3).Value = "\DataSet\NewUATDeals_DataSet\"
$ws.Cells.Item(6, 3).Value = "\DataSet\NewUATDeals_DataSet\"

# Widen the "UAT_Client" column so the longer client names fit.
$ws.Columns.Item(2).ColumnWidth = 16.6
